$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prevent Excel from auto-converting date-like text (column B) into date serials
$ws.Range("B3:B5").NumberFormat = "@"

# Row 3
$ws.Range("A3").Value2 = "SIT-00016"
$ws.Range("B3").Value2 = "2024-03-18"
$ws.Range("C3").Value2 = "09:36:59"
$ws.Range("D3").Value2 = "REVISE"
$ws.Range("E3").Value2 = "OK"
$ws.Range("F3").Value2 = "OK"
$ws.Range("G3").Value2 = "['HOME ?-29 : (114.009346, -7.691572000000011)']"
$ws.Range("H3").Value2 = "['HOME ?-1 : (114.010231, -7.692147000000005)', 'HOME ?-29 : (114.009346, -7.691572000000011)', 'HOME ?-49 : (114.006126, -7.691618)', 'HOME 5 : (114.006401, -7.692279000000015)', 'HOME 14 : (114.006065, -7.692483000000008)', 'HOME 13 : (114.006142, -7.69239099999999)', 'HOME ?-56 : (114.006233, -7.692200000000014)', 'HOME 9* : (114.006256, -7.692561000000007)', 'HOME 12 : (114.006157, -7.692311999999988)', 'HOME 11 : (114.006233, -7.692688000000003)', 'HOME ?-117 : (114.00631, -7.692485000000014)', 'HOME ?-146 : (114.008286, -7.69253)', 'HOME ?-150 : (114.010101, -7.692127999999991)', 'HOME ?-153 : (114.010132, -7.692245000000019)', 'HOME ?-163 : (114.009918, -7.692179000000015)', 'HOME ?-169 : (114.010099, -7.692218999999987)', 'HOME ?-170 : (114.009888, -7.692055000000007)', 'HOME ?-179 : (114.010033, -7.692185999999994)', 'HOME ?-180 : (114.008224, -7.692498000000009)', 'HOME ?-183 : (114.009857, -7.692186999999997)', 'HOME ?-186 : (114.009979, -7.692085999999994)', 'HOME ?-187 : (114.010155, -7.692352)', 'HOME 149 : (114.010132, -7.692450999999998)', 'HOME ?-207 : (114.009979, -7.692191000000008)']"
$ws.Range("I3").Value2 = "OK"
$ws.Range("J3").Value2 = "OK"
$ws.Range("K3").Value2 = "OK"
$ws.Range("L3").Value2 = "OK"
$ws.Range("M3").Value2 = "Revise"
$ws.Range("N3").Value2 = "OK"
$ws.Range("O3").Value2 = "Revise"
$ws.Range("P3").Value2 = "Revise"
$ws.Range("Q3").Value2 = "Revise"
$ws.Range("R3").Value2 = "Revise"
$ws.Range("S3").Value2 = "Revise"
$ws.Range("T3").Value2 = "Revise"
$ws.Range("U3").Value2 = "Revise"
$ws.Range("V3").Value2 = "Revise"
$ws.Range("W3").Value2 = "Revise"
$ws.Range("X3").Value2 = "Revise"
$ws.Range("Y3").Value2 = "Revise"
$ws.Range("Z3").Value2 = "Revise"
$ws.Range("AA3").Value2 = "Revise"
$ws.Range("AB3").Value2 = "Revise"
$ws.Range("AC3").Value2 = "Revise"
$ws.Range("AD3").Value2 = "Revise"
$ws.Range("AE3").Value2 = "Revise"
$ws.Range("AF3").Value2 = "Revise"
$ws.Range("AG3").Value2 = "OK"
$ws.Range("AH3").Value2 = "OK"
$ws.Range("AI3").Value2 = "OK"
$ws.Range("AJ3").Value2 = "OK"
$ws.Range("AK3").Value2 = "OK"
$ws.Range("AL3").Value2 = "OK"
$ws.Range("AM3").Value2 = "OK"
$ws.Range("AN3").Value2 = "OK"
$ws.Range("AO3").Value2 = "OK"
$ws.Range("AP3").Value2 = "OK"
$ws.Range("AQ3").Value2 = "OK"
$ws.Range("AR3").Value2 = "OK"
$ws.Range("AS3").Value2 = "OK"
$ws.Range("AT3").Value2 = "OK"
$ws.Range("AU3").Value2 = "OK"
$ws.Range("AV3").Value2 = "OK"
$ws.Range("AW3").Value2 = "OK"
$ws.Range("AX3").Value2 = "OK"
$ws.Range("AY3").Value2 = "OK"
$ws.Range("AZ3").Value2 = "OK"
$ws.Range("BA3").Value2 = "OK"
$ws.Range("BB3").Value2 = "OK"
$ws.Range("BC3").Value2 = "OK"
$ws.Range("BD3").Value2 = "OK"
$ws.Range("BE3").Value2 = "OK"
$ws.Range("BF3").Value2 = "OK"

# Row 4
$ws.Range("A4").Value2 = "SIT-00016"
$ws.Range("B4").Value2 = "2024-03-18"
$ws.Range("C4").Value2 = "09:39:55"
$ws.Range("D4").Value2 = "REVISE"
$ws.Range("E4").Value2 = "OK"
$ws.Range("F4").Value2 = "OK"
$ws.Range("G4").Value2 = "['HOME ?-29 : (114.009346, -7.691572000000011)']"
$ws.Range("H4").Value2 = "['HOME ?-1 : (114.010231, -7.692147000000005)', 'HOME ?-29 : (114.009346, -7.691572000000011)', 'HOME ?-49 : (114.006126, -7.691618)', 'HOME 5 : (114.006401, -7.692279000000015)', 'HOME 14 : (114.006065, -7.692483000000008)', 'HOME 13 : (114.006142, -7.69239099999999)', 'HOME ?-56 : (114.006233, -7.692200000000014)', 'HOME 9* : (114.006256, -7.692561000000007)', 'HOME 12 : (114.006157, -7.692311999999988)', 'HOME 11 : (114.006233, -7.692688000000003)', 'HOME ?-117 : (114.00631, -7.692485000000014)', 'HOME ?-146 : (114.008286, -7.69253)', 'HOME ?-150 : (114.010101, -7.692127999999991)', 'HOME ?-153 : (114.010132, -7.692245000000019)', 'HOME ?-163 : (114.009918, -7.692179000000015)', 'HOME ?-169 : (114.010099, -7.692218999999987)', 'HOME ?-170 : (114.009888, -7.692055000000007)', 'HOME ?-179 : (114.010033, -7.692185999999994)', 'HOME ?-180 : (114.008224, -7.692498000000009)', 'HOME ?-183 : (114.009857, -7.692186999999997)', 'HOME ?-186 : (114.009979, -7.692085999999994)', 'HOME ?-187 : (114.010155, -7.692352)', 'HOME 149 : (114.010132, -7.692450999999998)', 'HOME ?-207 : (114.009979, -7.692191000000008)']"
$ws.Range("I4").Value2 = "OK"
$ws.Range("J4").Value2 = "OK"
$ws.Range("K4").Value2 = "OK"
$ws.Range("L4").Value2 = "OK"
$ws.Range("M4").Value2 = "Revise"
$ws.Range("N4").Value2 = "OK"
$ws.Range("O4").Value2 = "Revise"
$ws.Range("P4").Value2 = "Revise"
$ws.Range("Q4").Value2 = "Revise"
$ws.Range("R4").Value2 = "Revise"
$ws.Range("S4").Value2 = "Revise"
$ws.Range("T4").Value2 = "Revise"
$ws.Range("U4").Value2 = "Revise"
$ws.Range("V4").Value2 = "Revise"
$ws.Range("W4").Value2 = "Revise"
$ws.Range("X4").Value2 = "Revise"
$ws.Range("Y4").Value2 = "Revise"
$ws.Range("Z4").Value2 = "Revise"
$ws.Range("AA4").Value2 = "Revise"
$ws.Range("AB4").Value2 = "Revise"
$ws.Range("AC4").Value2 = "Revise"
$ws.Range("AD4").Value2 = "Revise"
$ws.Range("AE4").Value2 = "Revise"
$ws.Range("AF4").Value2 = "Revise"
$ws.Range("AG4").Value2 = "OK"
$ws.Range("AH4").Value2 = "OK"
$ws.Range("AI4").Value2 = "OK"
$ws.Range("AJ4").Value2 = "OK"
$ws.Range("AK4").Value2 = "OK"
$ws.Range("AL4").Value2 = "OK"
$ws.Range("AM4").Value2 = "OK"
$ws.Range("AN4").Value2 = "OK"
$ws.Range("AO4").Value2 = "OK"
$ws.Range("AP4").Value2 = "OK"
$ws.Range("AQ4").Value2 = "OK"
$ws.Range("AR4").Value2 = "OK"
$ws.Range("AS4").Value2 = "OK"
$ws.Range("AT4").Value2 = "OK"
$ws.Range("AU4").Value2 = "OK"
$ws.Range("AV4").Value2 = "OK"
$ws.Range("AW4").Value2 = "OK"
$ws.Range("AX4").Value2 = "OK"
$ws.Range("AY4").Value2 = "OK"
$ws.Range("AZ4").Value2 = "OK"
$ws.Range("BA4").Value2 = "OK"
$ws.Range("BB4").Value2 = "OK"
$ws.Range("BC4").Value2 = "OK"
$ws.Range("BD4").Value2 = "OK"
$ws.Range("BE4").Value2 = "OK"
$ws.Range("BF4").Value2 = "OK"

# Row 5
$ws.Range("A5").Value2 = "SIT-00016"
$ws.Range("B5").Value2 = "2024-03-18"
$ws.Range("C5").Value2 = "09:49:55"
$ws.Range("D5").Value2 = "REVISE"
$ws.Range("E5").Value2 = "OK"
$ws.Range("F5").Value2 = "OK"
$ws.Range("G5").Value2 = "['HOME ?-29 : (114.009346, -7.691572000000011)']"
$ws.Range("H5").Value2 = "['HOME ?-1 : (114.010231, -7.692147000000005)', 'HOME ?-29 : (114.009346, -7.691572000000011)', 'HOME ?-49 : (114.006126, -7.691618)', 'HOME 5 : (114.006401, -7.692279000000015)', 'HOME 14 : (114.006065, -7.692483000000008)', 'HOME 13 : (114.006142, -7.69239099999999)', 'HOME ?-56 : (114.006233, -7.692200000000014)', 'HOME 9* : (114.006256, -7.692561000000007)', 'HOME 12 : (114.006157, -7.692311999999988)', 'HOME 11 : (114.006233, -7.692688000000003)', 'HOME ?-117 : (114.00631, -7.692485000000014)', 'HOME ?-146 : (114.008286, -7.69253)', 'HOME ?-150 : (114.010101, -7.692127999999991)', 'HOME ?-153 : (114.010132, -7.692245000000019)', 'HOME ?-163 : (114.009918, -7.692179000000015)', 'HOME ?-169 : (114.010099, -7.692218999999987)', 'HOME ?-170 : (114.009888, -7.692055000000007)', 'HOME ?-179 : (114.010033, -7.692185999999994)', 'HOME ?-180 : (114.008224, -7.692498000000009)', 'HOME ?-183 : (114.009857, -7.692186999999997)', 'HOME ?-186 : (114.009979, -7.692085999999994)', 'HOME ?-187 : (114.010155, -7.692352)', 'HOME 149 : (114.010132, -7.692450999999998)', 'HOME ?-207 : (114.009979, -7.692191000000008)']"
$ws.Range("I5").Value2 = "OK"
$ws.Range("J5").Value2 = "OK"
$ws.Range("K5").Value2 = "OK"
$ws.Range("L5").Value2 = "OK"
$ws.Range("M5").Value2 = "Revise"
$ws.Range("N5").Value2 = "OK"
$ws.Range("O5").Value2 = "Revise"
$ws.Range("P5").Value2 = "Revise"
$ws.Range("Q5").Value2 = "Revise"
$ws.Range("R5").Value2 = "Revise"
$ws.Range("S5").Value2 = "Revise"
$ws.Range("T5").Value2 = "Revise"
$ws.Range("U5").Value2 = "Revise"
$ws.Range("V5").Value2 = "Revise"
$ws.Range("W5").Value2 = "Revise"
$ws.Range("X5").Value2 = "Revise"
$ws.Range("Y5").Value2 = "Revise"
$ws.Range("Z5").Value2 = "Revise"
$ws.Range("AA5").Value2 = "Revise"
$ws.Range("AB5").Value2 = "Revise"
$ws.Range("AC5").Value2 = "Revise"
$ws.Range("AD5").Value2 = "Revise"
$ws.Range("AE5").Value2 = "Revise"
$ws.Range("AF5").Value2 = "Revise"
$ws.Range("AG5").Value2 = "OK"
$ws.Range("AH5").Value2 = "OK"
$ws.Range("AI5").Value2 = "OK"
$ws.Range("AJ5").Value2 = "OK"
$ws.Range("AK5").Value2 = "OK"
$ws.Range("AL5").Value2 = "OK"
$ws.Range("AM5").Value2 = "OK"
$ws.Range("AN5").Value2 = "OK"
$ws.Range("AO5").Value2 = "OK"
$ws.Range("AP5").Value2 = "OK"
$ws.Range("AQ5").Value2 = "OK"
$ws.Range("AR5").Value2 = "OK"
$ws.Range("AS5").Value2 = "OK"
$ws.Range("AT5").Value2 = "OK"
$ws.Range("AU5").Value2 = "OK"
$ws.Range("AV5").Value2 = "OK"
$ws.Range("AW5").Value2 = "OK"
$ws.Range("AX5").Value2 = "OK"
$ws.Range("AY5").Value2 = "OK"
$ws.Range("AZ5").Value2 = "OK"
$ws.Range("BA5").Value2 = "OK"
$ws.Range("BB5").Value2 = "OK"
$ws.Range("BC5").Value2 = "OK"
$ws.Range("BD5").Value2 = "OK"
$ws.Range("BE5").Value2 = "OK"
$ws.Range("BF5").Value2 = "OK"

# Reset number format on column B back to default so no style residue remains
$ws.Range("B3:B5").Style = "Normal"
